# homework5.docx update:
# ". The graph represents observed attendance at 14 social events by 18
#  Southern "high society" women from the South."
# becomes a run-split sentence with "14" -> "8" and "18" -> "10":
# ". The graph represents observed attendance at " / "8" / " social events by " /
# "10" / " Southern "high society" women from the South."

$d = $word.ActiveDocument

# Curly quotes used around "high society" in the source text.
$lq = [char]0x201C
$rq = [char]0x201D

# Locate the sentence fragment that needs to change (from the first number
# through the end of the paragraph) using Find, the same way a user would.
$old = "14 social events by 18 Southern " + $lq + "high society" + $rq + " women from the South."
$hit = $d.Content
$found = $hit.Find.Execute($old)
if (-not $found) {
    throw "Could not locate target sentence to edit."
}

# Re-seat a plain Range over the exact hit so the edit is anchored to real
# character offsets (not a live Find object).
$editStart = $hit.Start
$editEnd = $hit.End
$target = $d.Range($editStart, $editEnd)

# Build the five replacement runs (identical, default run formatting, just
# like the original single run) that together replace the old wording.
$newRunsXml = `
    '<w:r><w:t xml:space="preserve">8</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> social events by </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">10</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> Southern ' + $lq + 'high society' + $rq + ' women from the South.</w:t></w:r>'

$packageXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $newRunsXml + '</w:p>'

$target.InsertXML($packageXml)

Write-Output "Updated attendance sentence (14 -> 8 events, 18 -> 10 women)."
